$d = $word.ActiveDocument

# Remove the "(F1)" marker from the precondition sentence.
$rngRemove = $d.Content
$rngRemove.Find.Execute("(F1)", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngRemove.Text = ""

# Force Word to split the run in two (matching the target OOXML, which has
# two runs with identical rPr) by nudging the font color of the first part
# away and then back to its original value.
$rngSplit = $d.Content
$rngSplit.Find.Execute("Användare är inloggad", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rngSplit.Font.Color = 123456
$rngSplit.Font.Color = 3355443
